{"js": "// Fill in the first empty timesheet row (right after the 25/11/2012 \"2pm-5pm\"\n// row) with the new entry: date, hours, and the \"what done\" text \u2014 the latter\n// keeps the spell-check markers + the relocated _GoBack bookmark + trailing\n// space run, matching the authored OOXML.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row 0 = header. Find the first fully-blank row after that (it comes right\n// after the existing \"25/11/2012 / 2pm-5pm / Lots of...\" row) to populate.\nlet targetIndex = -1;\nfor (let i = 1; i < table.values.length; i++) {\n  if (table.values[i].every((c) => c.trim() === \"\")) {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error(\"No blank timesheet row found to fill in\");\n}\n\nconst targetRow = rows.items[targetIndex];\nconst cells = targetRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\nconst dateCell = cells.items[0];\nconst hoursCell = cells.items[1];\nconst whatDoneCell = cells.items[2];\n\ndateCell.body.insertText(\"25/11/2012\", Word.InsertLocation.replace);\nhoursCell.body.insertText(\"6pm-7:30pm\", Word.InsertLocation.replace);\n\n// The \"what done\" cell needs the proofErr spell-check wrapper around\n// \"friend,battle,breed\" plus the _GoBack bookmark before the trailing space,\n// so build it with raw OOXML instead of plain text.\nconst whatDoneOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">Added ability to make </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>friend,battle,breed</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> requests</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nwhatDoneCell.body.insertOoxml(whatDoneOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Fill in the first empty timesheet row (right after the 25/11/2012\n# \"2pm-5pm\" row) with the new entry: date, hours, and the \"what done\" text.\n# The \"what done\" cell needs the proofErr spell-check wrapper around\n# \"friend,battle,breed\" plus the relocated _GoBack bookmark before a\n# trailing space run, so it is written via raw OOXML (InsertXML) rather\n# than plain Range.Text, and \u2014 for a byte-for-byte match with the other\n# two cells \u2014 they use the same InsertXML path too.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Locate the first fully-blank data row after the header/first data row.\n$targetRow = 0\nfor ($i = 2; $i -le $t.Rows.Count; $i++) {\n    $c1 = ($t.Cell($i, 1).Range.Text) -replace \"[\\x07\\x0d\\x0c]\", \"\"\n    $c2 = ($t.Cell($i, 2).Range.Text) -replace \"[\\x07\\x0d\\x0c]\", \"\"\n    $c3 = ($t.Cell($i, 3).Range.Text) -replace \"[\\x07\\x0d\\x0c]\", \"\"\n    if ($c1.Trim() -eq \"\" -and $c2.Trim() -eq \"\" -and $c3.Trim() -eq \"\") {\n        $targetRow = $i\n        break\n    }\n}\n\nif ($targetRow -eq 0) {\n    throw \"No blank timesheet row found to fill in\"\n}\n\nfunction Get-CellPackageXml($innerBodyXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body>' + $innerBodyXml + '</w:body>' +\n      '</w:document>' +\n      '</pkg:xmlData>' +\n      '</pkg:part>' +\n      '</pkg:package>'\n}\n\n$dateXml = Get-CellPackageXml '<w:p><w:r><w:t>25/11/2012</w:t></w:r></w:p>'\n[void]$t.Cell($targetRow, 1).Range.InsertXML($dateXml)\n\n$hoursXml = Get-CellPackageXml '<w:p><w:r><w:t>6pm-7:30pm</w:t></w:r></w:p>'\n[void]$t.Cell($targetRow, 2).Range.InsertXML($hoursXml)\n\n$whatDoneBody = '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">Added ability to make </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>friend,battle,breed</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> requests</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '</w:p>'\n$whatDoneXml = Get-CellPackageXml $whatDoneBody\n[void]$t.Cell($targetRow, 3).Range.InsertXML($whatDoneXml)\n"}
